$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 398.4
$ws.Range("J28").Value = 998
$ws.Range("L28").Value = 998
$ws.Range("N28").Value = -1968

$ws.Range("H53").Value = 7896.6665
$ws.Range("I53").Value = 9754.083000000001
$ws.Range("K53").Value = 9754.083000000001
$ws.Range("M53").Value = -9117.083000000001

$ws.Range("H88").Value = 4129.857
$ws.Range("I88").Value = 2005
$ws.Range("J88").Value = 4979.8
$ws.Range("K88").Value = 2005
$ws.Range("L88").Value = 4979.8
$ws.Range("M88").Value = -1599
$ws.Range("N88").Value = -5791.8

$ws.Range("H91").Value = 4129.857
$ws.Range("I91").Value = 2005
$ws.Range("J91").Value = 4979.8
$ws.Range("K91").Value = 2005
$ws.Range("L91").Value = 4979.8
$ws.Range("M91").Value = -601
$ws.Range("N91").Value = -7787.8

$ws.Range("H107").Value = 761.5833
$ws.Range("I107").Value = 477.14285
$ws.Range("K107").Value = 477.14285
$ws.Range("M107").Value = 1442.85715

$ws.Range("H113").Value = 36666.668
$ws.Range("J113").Value = 3500
$ws.Range("L113").Value = 3500
$ws.Range("N113").Value = -10008

$ws.Range("H137").Value = 2024.9375
$ws.Range("I137").Value = 1249.7693
$ws.Range("J137").Value = 2555.3157
$ws.Range("K137").Value = 3749.3079
$ws.Range("L137").Value = 7665.9471
$ws.Range("M137").Value = -1199.3079
$ws.Range("N137").Value = -12765.9471

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3322915.8
$ws.Range("J2").Value = 833
$ws.Range("L2").Value = 833
$ws.Range("N2").Value = -1059

$ws.Range("H32").Value = 3942.75
$ws.Range("I32").Value = 3450.878
$ws.Range("J32").Value = 10665
$ws.Range("K32").Value = 3450.878
$ws.Range("L32").Value = 10665
$ws.Range("M32").Value = -3163.878
$ws.Range("N32").Value = -11239

$ws.Range("H74").Value = 1455.3334
$ws.Range("I74").Value = 1436
$ws.Range("J74").Value = 1479.5
$ws.Range("K74").Value = 1436
$ws.Range("L74").Value = 1479.5
$ws.Range("M74").Value = -562
$ws.Range("N74").Value = -3227.5

$ws.Range("H77").Value = 1455.3334
$ws.Range("I77").Value = 1436
$ws.Range("J77").Value = 1479.5
$ws.Range("K77").Value = 7180
$ws.Range("L77").Value = 7397.5
$ws.Range("M77").Value = -2812
$ws.Range("N77").Value = -16133.5

$ws.Range("H116").Value = 3322915.8
$ws.Range("J116").Value = 833
$ws.Range("L116").Value = 833
$ws.Range("N116").Value = -5421

$ws.Range("H132").Value = 1280.965
$ws.Range("I132").Value = 783.93616
$ws.Range("J132").Value = 3617
$ws.Range("K132").Value = 2351.80848
$ws.Range("L132").Value = 10851
$ws.Range("M132").Value = 178.1915200000003
$ws.Range("N132").Value = -15911

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3322915.8
$ws.Range("J3").Value = 833
$ws.Range("L3").Value = 833
$ws.Range("N3").Value = -1061

$ws.Range("H20").Value = 1647.129
$ws.Range("I20").Value = 1661.875
$ws.Range("J20").Value = 1596.5714
$ws.Range("K20").Value = 1661.875
$ws.Range("L20").Value = 1596.5714
$ws.Range("M20").Value = -1414.875
$ws.Range("N20").Value = -2090.5714

$ws.Range("H105").Value = 2608.3157
$ws.Range("I105").Value = 2468.5715
$ws.Range("K105").Value = 2468.5715
$ws.Range("M105").Value = -721.5715

$ws.Range("H134").Value = 4292.615
$ws.Range("I134").Value = 4687.024
$ws.Range("J134").Value = 2636.1
$ws.Range("K134").Value = 14061.072
$ws.Range("L134").Value = 7908.299999999999
$ws.Range("M134").Value = -11526.072
$ws.Range("N134").Value = -12978.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 320
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 466.66666
$ws.Range("K7").Value = 100
$ws.Range("L7").Value = 466.66666
$ws.Range("M7").Value = 13
$ws.Range("N7").Value = -692.66666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 8834.5
$ws.Range("I56").Value = 8834.5
$ws.Range("K56").Value = 8834.5
$ws.Range("M56").Value = -8304.5

$ws.Range("H68").Value = 721.0909
$ws.Range("I68").Value = 716
$ws.Range("J68").Value = 725.3333
$ws.Range("K68").Value = 2148
$ws.Range("L68").Value = 2175.9999
$ws.Range("M68").Value = -1337
$ws.Range("N68").Value = -3797.9999

$ws.Range("H71").Value = 721.0909
$ws.Range("I71").Value = 716
$ws.Range("J71").Value = 725.3333
$ws.Range("K71").Value = 6444
$ws.Range("L71").Value = 6527.9997
$ws.Range("M71").Value = -2388
$ws.Range("N71").Value = -14639.9997

$ws.Range("H98").Value = 661.375
$ws.Range("I98").Value = 248.66667
$ws.Range("J98").Value = 909
$ws.Range("K98").Value = 746.00001
$ws.Range("L98").Value = 2727
$ws.Range("M98").Value = 751.99999
$ws.Range("N98").Value = -5723

$ws.Range("H109").Value = 22730972
$ws.Range("I109").Value = 55556616
$ws.Range("J109").Value = 5524.769
$ws.Range("K109").Value = 166669848
$ws.Range("L109").Value = 16574.307
$ws.Range("M109").Value = -166668808
$ws.Range("N109").Value = -18654.307

$ws.Range("H131").Value = 7587204.5
$ws.Range("J131").Value = 12530.3
$ws.Range("L131").Value = 37590.89999999999
$ws.Range("N131").Value = -47670.89999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2096.8096
$ws.Range("I102").Value = 2090.853
$ws.Range("J102").Value = 2122.125
$ws.Range("K102").Value = 2090.853
$ws.Range("L102").Value = 2122.125
$ws.Range("M102").Value = -468.8530000000001
$ws.Range("N102").Value = -5366.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 347.1154
$ws.Range("I55").Value = 303.31818
$ws.Range("J55").Value = 588
$ws.Range("K55").Value = 303.31818
$ws.Range("L55").Value = 588
$ws.Range("M55").Value = -130.31818
$ws.Range("N55").Value = -934

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 747.6923
$ws.Range("I113").Value = 442
$ws.Range("J113").Value = 1766.6666
$ws.Range("K113").Value = 1326
$ws.Range("L113").Value = 5299.9998
$ws.Range("M113").Value = 844
$ws.Range("N113").Value = -9639.9998

$ws.Range("H132").Value = 1642.6444
$ws.Range("I132").Value = 993.3
$ws.Range("J132").Value = 2941.3333
$ws.Range("K132").Value = 2979.9
$ws.Range("L132").Value = 8823.999899999999
$ws.Range("M132").Value = -449.8999999999996
$ws.Range("N132").Value = -13883.9999

